$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header F1: "Потерято" -> "Потеряно"
$ws.Range("F1").Value = "Потеряно, байт"

# Row 2 numbers updated to reflect the integrated test duration in express-test.
# C2 is stored as text (t="str") in the original file, so write it through a
# text-producing formula and then collapse that formula down to its cached
# literal value (Copy + PasteSpecial values) - this keeps the cell a plain
# text value (no formula, no stray number-format/style change) instead of
# Excel's usual "value looks numeric -> coerce to Number" behavior.
$ws.Range("C2").Formula = '="26.899999618530273"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("D2").Value = 115416192
$ws.Range("E2").Value = 113735808
$ws.Range("F2").Value = 1680384
$ws.Range("G2").Value = 1.4559343631784352
